$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the two stray <w:lastRenderedPageBreak/> markers that precede
#    "Type of Data the model was given:" and the "F" (of "Few-shot learning
#    and Finetuning:") runs. A Find/Replace that rewrites the run's text
#    (even to the same text) regenerates the run without the stale
#    last-rendered-page-break marker while keeping the run formatting.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Type of Data the model was given:", $true, $false, $false, $false, $false, $true, 1, $false, "Type of Data the model was given:", 2) | Out-Null

# The single-letter "F" run needs to be targeted precisely (searching for
# just "F" anywhere would be ambiguous), so locate the paragraph that
# contains "Few-shot learning and Finetuning:" and restrict the Find to its
# very first character.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Few-shot learning and Finetuning*") {
        $firstChar = $d.Range($para.Range.Start, $para.Range.Start + 1)
        $firstChar.Find.Execute("F", $true, $false, $false, $false, $false, $true, 1, $false, "F", 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Add <w:noProof/> to the run that carries the last image's
#    lastRenderedPageBreak + drawing (the picture immediately followed by
#    the "Snippet from Database" caption). All the other inline pictures in
#    the document already have NoProofing set; this is the one missing it.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -like "*Snippet from Database*") {
            $para.Range.NoProofing = $true
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Insert two new paragraphs right after "I have discussed about my
#    current approach..." (and before the pre-existing pair of empty
#    paragraphs): a blank paragraph, then one holding the new sentence
#    "My application works both as me and can respond to general questions
#    also like chatGPT".
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*I have discussed about my current approach*") {
        $para.Range.InsertParagraphAfter()
        $blankPara = $d.Paragraphs.Item($i + 1)
        $blankPara.Range.InsertParagraphAfter()
        $textPara = $d.Paragraphs.Item($i + 2)
        $textPara.Range.Text = "My application works both as me and can respond to general questions also like chatGPT"
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Append a new sentence at the very end of the document, right after
#    "...the limitation to the API access, I could do this."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("the limitation to the API access, I could do this.", $true, $false, $false, $false, $false, $true, 1, $false, "the limitation to the API access, I could do this. But not only this , further the model can be finetuned or few shot trained on different models with good amount of resources available.", 2) | Out-Null

"done"
